$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-5 from 45208 to 45212
$ws.Range("C2").Value2 = 45212
$ws.Range("C3").Value2 = 45212
$ws.Range("C4").Value2 = 45212
$ws.Range("C5").Value2 = 45212
